# feat: add 2022-Q1 data
#
# The workbook currently has sheets: 2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3,
# 2021-Q4, 总计 (Total). We need to insert a new "2022-Q1" sheet with the
# quarterly fund-holdings detail (mirroring the format of "2021-Q4"), right
# before the "总计" sheet, and update "总计" with a new summary row for
# 2022-Q1 (pushing the existing summary rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: re-create the "总计" sheet so that the new "2022-Q1" sheet ends
# up with sheetId 6 (re-using the id vacated by 总计) and "总计" itself
# gets sheetId 7, matching how Excel allocates ids when a sheet is
# removed and two fresh sheets are appended afterwards.
# ---------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete()

$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)

$q1Sheet = $wb.Worksheets.Add($null, $lastSheet)
$q1Sheet.Name = "2022-Q1"

$totalSheet = $wb.Worksheets.Add($null, $q1Sheet)
$totalSheet.Name = "总计"

# Sheet used to source the header/number formatting for the fund table
$refSheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# Step 2: populate the "2022-Q1" sheet with the fund-holdings detail.
# ---------------------------------------------------------------------
$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Copy header/column formatting from the equivalent row on 2021-Q4 so the
# new cells share the same style (bold, centered, bordered) as the rest
# of the workbook instead of Excel's plain default.
$refSheet.Range("B1:H1").Copy()
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)

$q1Data = @(
    @(0,  "005888", "华夏新兴消费混合A",              "16.37", "91.96", "5.36", "0.8774", 2),
    @(1,  "519091", "新华泛资源优势混合",              "13.39", "78.42", "5.07", "0.6789", 5),
    @(2,  "009885", "新华景气行业混合A",                "9.84", "85.63", "5.01", "0.4930", 5),
    @(3,  "519089", "新华优选成长混合",                  "8.58", "87.20", "5.43", "0.4659", 4),
    @(4,  "014150", "新华鑫益灵活配置混合A",             "7.33", "84.57", "5.09", "0.3731", 5),
    @(5,  "000584", "新华鑫益灵活配置混合",              "5.89", "84.57", "5.09", "0.2998", 5),
    @(6,  "005889", "华夏新兴消费混合C",                 "4.25", "91.96", "5.36", "0.2278", 2),
    @(7,  "009886", "新华景气行业混合C",                 "1.28", "85.63", "5.01", "0.0641", 5),
    @(8,  "009829", "长城优选增强六个月持有期混合A",      "5.36", "24.94", "0.95", "0.0509", 8),
    @(9,  "010799", "长城优选稳进六个月持有期混合A",      "2.70", "37.25", "1.54", "0.0416", 8),
    @(10, "011538", "长城优选添瑞六个月持有期混合A",      "4.43", "25.71", "0.93", "0.0412", 8),
    @(11, "004573", "新华鑫泰灵活配置混合",               "0.73", "77.81", "4.93", "0.0360", 3),
    @(12, "009169", "湘财长兴灵活配置混合A",              "1.16", "85.40", "2.85", "0.0331", 10),
    @(13, "005910", "广发龙头优选灵活配置混合",           "0.50", "92.44", "5.02", "0.0251", 7),
    @(14, "009170", "湘财长兴灵活配置混合C",              "0.46", "85.40", "2.85", "0.0131", 10),
    @(15, "004189", "华商消费行业股票",                   "0.30", "81.87", "4.22", "0.0127", 5),
    @(16, "002543", "长城久益灵活配置混合A",              "0.36", "89.88", "3.40", "0.0122", 10),
    @(17, "165524", "信诚中证智能家居指数（LOF）",        "0.40", "93.89", "1.40", "0.0056", 2),
    @(18, "002544", "长城久益灵活配置混合C",              "0.07", "89.88", "3.40", "0.0024", 10),
    @(19, "009830", "长城优选增强六个月持有期混合C",      "0.25", "24.94", "0.95", "0.0024", 8),
    @(20, "010800", "长城优选稳进六个月持有期混合C",      "0.06", "37.25", "1.54", "0.0009", 8),
    @(21, "673120", "西部利得新富灵活配置混合",           "0.02", "80.77", "4.09", "0.0008", 7),
    @(22, "011539", "长城优选添瑞六个月持有期混合C",      "0.01", "25.71", "0.93", "0.0001", 8)
)

$row = 2
foreach ($item in $q1Data) {
    $q1Sheet.Range("A$row").Value = $item[0]
    $q1Sheet.Range("B$row").Value = "'" + $item[1]
    $q1Sheet.Range("C$row").Value = $item[2]
    $q1Sheet.Range("D$row").Value = "'" + $item[3]
    $q1Sheet.Range("E$row").Value = "'" + $item[4]
    $q1Sheet.Range("F$row").Value = "'" + $item[5]
    $q1Sheet.Range("G$row").Value = "'" + $item[6]
    $q1Sheet.Range("H$row").Value = $item[7]
    $row = $row + 1
}

# Apply the "index column" style (bold / centered / bordered, matching A2
# on 2021-Q4) to every A-column cell of the new data rows.
$refSheet.Range("A2").Copy()
$q1Sheet.Range("A2:A24").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 3: rebuild the "总计" summary sheet - add the new 2022-Q1 row on
# top and push the previously existing rows down by one.
# ---------------------------------------------------------------------
$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$refTotalHeader = $refSheet.Range("B1:D1")
$refTotalHeader.Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122)

$totalData = @(
    @(0, "2022-Q1", 23, 3.76),
    @(1, "2021-Q4", 18, 6.11),
    @(2, "2021-Q3", 5,  0.89),
    @(3, "2021-Q2", 18, 2.31),
    @(4, "2021-Q1", 32, 5.79),
    @(5, "2020-Q4", 8,  0.9)
)

$row = 2
foreach ($item in $totalData) {
    $totalSheet.Range("A$row").Value = $item[0]
    $totalSheet.Range("B$row").Value = $item[1]
    $totalSheet.Range("C$row").Value = $item[2]
    $totalSheet.Range("D$row").Value = $item[3]
    $row = $row + 1
}

$refSheet.Range("A2").Copy()
$totalSheet.Range("A2:A7").PasteSpecial(-4122)
